# Nowcasts 2025Q4 update
# The table on Sheet1 rolls forward to a new set of six "vintage" dates
# (2025-09-30 .. 2025-12-15). Row 1 keeps the same header labels; rows 2-7
# get the new dates (column A) plus refreshed Prognose/Revision figures
# (columns B-K). Rows 8-11 are left exactly as they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header labels (unchanged text, rewritten so the shared-string table
# regenerates cleanly alongside the new rows below)
$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "Prognose"
$ws.Range("C1").Value = "surveys"
$ws.Range("D1").Value = "production"
$ws.Range("E1").Value = "orders"
$ws.Range("F1").Value = "turnover"
$ws.Range("G1").Value = "financial"
$ws.Range("H1").Value = "labor market"
$ws.Range("I1").Value = "prices"
$ws.Range("J1").Value = "national accounts"
$ws.Range("K1").Value = "Revision"

# Rows 2-7: new 2025Q4 vintage dates + nowcast/revision data
# (Column A holds the vintage date as plain text, exactly like the rest of
# the sheet - force text format first so Excel doesn't auto-convert the
# ISO-looking string into a real date serial, then drop back to the
# Normal style so no stray number-format survives on the cell.)
$ws.Range("A2:A7").NumberFormat = "@"

$ws.Range("A2").Value = "2025-09-30"
$ws.Range("B2").Value = 0.28065104601029489
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

$ws.Range("A3").Value = "2025-10-15"
$ws.Range("B3").Value = 0.29874339719721282
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.0026275856430242875
$ws.Range("E3").Value = 0.0014258307813590652
$ws.Range("F3").Value = 0.0013105236324631986
$ws.Range("G3").Value = 0.00098263665740082589
$ws.Range("H3").Value = -0.000092775578417227288
$ws.Range("I3").Value = -0.00048999626027333753
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.00047395028266944417

$ws.Range("A4").Value = "2025-10-30"
$ws.Range("B4").Value = 0.38479727604841291
$ws.Range("C4").Value = 0.053889086734521557
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = -0.00004629444607496037
$ws.Range("F4").Value = -0.000014127388037332344
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.00046820258468735357
$ws.Range("I4").Value = -0.0037477176633198021
$ws.Range("J4").Value = 0.0031744241999277566
$ws.Range("K4").Value = -0.0020280586086206531

$ws.Range("A5").Value = "2025-11-15"
$ws.Range("B5").Value = 0.37259765339111584
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = -0.00066439599030669475
$ws.Range("E5").Value = -0.00024837836398691304
$ws.Range("F5").Value = 0.0097909129909326507
$ws.Range("G5").Value = -0.0019837507965679469
$ws.Range("H5").Value = 0.00060866248784145489
$ws.Range("I5").Value = -0.00074744182122994077
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -0.019197580537502534

$ws.Range("A6").Value = "2025-11-30"
$ws.Range("B6").Value = 0.27181800583034171
$ws.Range("C6").Value = -0.030999586667594219
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = -0.0010196349353219615
$ws.Range("F6").Value = 0.00059424081433660022
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0.0011609445355531627
$ws.Range("I6").Value = -0.010549455370754916
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0.00002730309716314494

$ws.Range("A7").Value = "2025-12-15"
$ws.Range("B7").Value = 0.25266098949194038
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = -0.030496042536245997
$ws.Range("E7").Value = -0.0014820240562763538
$ws.Range("F7").Value = 0.007463536253667951
$ws.Range("G7").Value = 0.0020777276731186963
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = -0.0026969537486006834

$ws.Range("A2:A7").Style = "Normal"

# Rows 8-11 are intentionally left untouched (same vintage figures as before).

# Column width tweaks (auto-fit drift from the refreshed figures):
#   D: 15.24609375 -> 16.24609375, F: 16.24609375 -> 15.77734375, H: 16.24609375 -> 15.77734375
$ws.Columns.Item(4).ColumnWidth = 15.25
$ws.Columns.Item(6).ColumnWidth = 14.92
$ws.Columns.Item(8).ColumnWidth = 14.92
